# Append year-index data rows 199:210 (A=197..208, B=corresponding values)
# to the "월_수출물가지수" sheet, extending the used range from A1:B198 to A1:B210.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (A value, B value) pairs for new rows, in order starting at row 199
$data = @(
  @(197, 0.5264750378214825),
  @(198, 0.4208499518635676),
  @(199, 0.3548342731398707),
  @(200, 0.5998257919589235),
  @(201, 0.5099711181405583),
  @(202, 0.426508438611313),
  @(203, 0.5264750378214825),
  @(204, 0.4776234355659468),
  @(205, 0.4753128868106175),
  @(206, 0.576536927520286),
  @(207, 0.473662494842525),
  @(208, 0.473662494842525)
)

# Copy the formatting of the last existing data cell (A198: bold, bordered,
# centered style) so the new column-A cells keep the same look.
$ws.Range("A198").Copy()

$startRow = 199
for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $startRow + $i
  $ws.Range("A$row").PasteSpecial(-4122)
  $ws.Range("A$row").Value = $data[$i][0]
  $ws.Range("B$row").Value = $data[$i][1]
}

$excel.CutCopyMode = 0
